$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is inserted at row 89; the existing rows
# 89-94 shift down to 90-95 (their data is unchanged, only their row
# position moves).
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with the new record.
$ws.Cells.Item(89, 1).Value = 1
$ws.Cells.Item(89, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(89, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(89, 4).Value = 45166
$ws.Cells.Item(89, 5).Value = 15
$ws.Cells.Item(89, 6).Value = 100112031
$ws.Cells.Item(89, 7).Value = "Poroto verde"
$ws.Cells.Item(89, 8).Value = "Sin especificar"
$ws.Cells.Item(89, 9).Value = "Primera"
$ws.Cells.Item(89, 10).Value = 2400
$ws.Cells.Item(89, 11).Value = 1000
$ws.Cells.Item(89, 12).Value = 1200
$ws.Cells.Item(89, 13).Value = 1125
$ws.Cells.Item(89, 14).Value = "`$/kilo"
$ws.Cells.Item(89, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(89, 16).Value = 1125
$ws.Cells.Item(89, 17).Value = 1
$ws.Cells.Item(89, 18).Value = "Hortaliza"
